# Apply label corrections to column A of the active sheet (Sheet1 of lbl6.xlsx).
# Each row below flips a single binary label value (0 <-> 1) as captured by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = 0
$ws.Range("A87").Value = 1
$ws.Range("A112").Value = 0
$ws.Range("A116").Value = 0
$ws.Range("A120").Value = 0
$ws.Range("A185").Value = 0
$ws.Range("A269").Value = 0
$ws.Range("A322").Value = 0
$ws.Range("A324").Value = 1
$ws.Range("A509").Value = 1
$ws.Range("A516").Value = 1
$ws.Range("A749").Value = 1
$ws.Range("A784").Value = 0
$ws.Range("A854").Value = 1
$ws.Range("A884").Value = 0
$ws.Range("A928").Value = 0
$ws.Range("A973").Value = 0
$ws.Range("A1022").Value = 1
$ws.Range("A1024").Value = 0
$ws.Range("A1034").Value = 0
$ws.Range("A1080").Value = 0
$ws.Range("A1087:A1099").Value = 0
$ws.Range("A1102").Value = 1
$ws.Range("A1141").Value = 1
$ws.Range("A1156").Value = 1
$ws.Range("A1158").Value = 1
$ws.Range("A1164").Value = 0
$ws.Range("A1179").Value = 1
$ws.Range("A1181").Value = 1
$ws.Range("A1197:A1198").Value = 1
$ws.Range("A1201").Value = 0
$ws.Range("A1218").Value = 0
$ws.Range("A1223").Value = 0
$ws.Range("A1229").Value = 0
$ws.Range("A1234").Value = 0
$ws.Range("A1236").Value = 1
$ws.Range("A1261:A1265").Value = 0
$ws.Range("A1452").Value = 0
$ws.Range("A1602").Value = 1
